# Reorder the rows of the "classFields" sheet so that the field list for
# pl.piomin.base.domain.Order matches the new field/type order produced by
# implementing the standard relationship between microservices.
#
# Original row order (rows 2-8): source, status, customerId, productCount, productId, price, id
# New row order      (rows 2-8): status, customerId, price, id, source, productId, productCount
#
# Only columns B (Field Name) and D (Field Type) change; column A (Class Name)
# and column C (Field Modifier = "private") stay identical on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# New field name / field type values for rows 2 through 8.
$fields = @(
    @{ Name = "status";       Type = "java.lang.String" },
    @{ Name = "customerId";   Type = "java.lang.Long" },
    @{ Name = "price";        Type = "int" },
    @{ Name = "id";           Type = "java.lang.Long" },
    @{ Name = "source";       Type = "java.lang.String" },
    @{ Name = "productId";    Type = "java.lang.Long" },
    @{ Name = "productCount"; Type = "int" }
)

$row = 2
foreach ($f in $fields) {
    $ws.Cells.Item($row, 2).Value = $f.Name
    $ws.Cells.Item($row, 4).Value = $f.Type
    $row++
}
